# Updated cryptos list (refresh of Price/Volume(1h) columns), mirroring
# the scheduled GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting (avoid numeric auto-conversion) for the Price column
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.040.74'
$ws.Range('E2').Value = '  -1.05%  '
$ws.Range('D3').Value = '1.825.22'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.36%  '
$ws.Range('D5').Value = '311.58'
$ws.Range('E5').Value = '  -1.21%  '
$ws.Range('E6').Value = '  -0.35%  '
$ws.Range('D7').Value = '0.4401'
$ws.Range('E7').Value = '  +2.59%  '
$ws.Range('D8').Value = '0.3680'
$ws.Range('E8').Value = '  -0.40%  '
$ws.Range('D9').Value = '0.07271'
$ws.Range('E9').Value = '  +0.40%  '
$ws.Range('D10').Value = '0.8443'
$ws.Range('E10').Value = '  -2.27%  '
$ws.Range('D11').Value = '20.69'
$ws.Range('E11').Value = '  -2.03%  '
$ws.Range('D12').Value = '1.821.63'
$ws.Range('E12').Value = '  -0.10%  '
$ws.Range('E13').Value = '  -0.08%  '
$ws.Range('D14').Value = '0.07073'
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('D15').Value = '5.304'
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('D16').Value = '90.03'
$ws.Range('E16').Value = '  +2.85%  '
$ws.Range('D17').Value = '1.002'
$ws.Range('E17').Value = '  -0.38%  '
$ws.Range('D18').Value = '0.000008797'
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('D19').Value = '1.000'
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('E20').Value = '  -1.83%  '
$ws.Range('D21').Value = '27.018.49'
$ws.Range('E21').Value = '  -1.23%  '
$ws.Range('D22').Value = '5.151'
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').Value = '10.89'
$ws.Range('E23').Value = '  +0.37%  '
$ws.Range('D24').Value = '2.044.14'
$ws.Range('E24').Value = '  -0.32%  '
$ws.Range('D25').Value = '1.993'
$ws.Range('E25').Value = '  -0.92%  '
$ws.Range('E26').Value = '  -0.88%  '
$ws.Range('D27').Value = '2.203'
$ws.Range('E27').Value = '  +3.04%  '
$ws.Range('E28').Value = '  -0.57%  '
$ws.Range('D29').Value = '5.231'
$ws.Range('E29').Value = '  -0.82%  '
$ws.Range('D30').Value = '116.97'
$ws.Range('E30').Value = '  +0.26%  '
$ws.Range('D31').Value = '0.08797'
$ws.Range('E31').Value = '  -0.64%  '
$ws.Range('D32').Value = '1.180'
$ws.Range('E32').Value = '  -1.61%  '
$ws.Range('D33').Value = '0.7405'
$ws.Range('E33').Value = '  -3.00%  '
$ws.Range('D34').Value = '4.425'
$ws.Range('E34').Value = '  -1.31%  '
$ws.Range('D35').Value = '2.885'
$ws.Range('E35').Value = '  +1.14%  '
$ws.Range('D36').Value = '0.9999'
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('E37').Value = '  -2.19%  '
$ws.Range('E38').Value = '  -0.39%  '
$ws.Range('D39').Value = '0.05236'
$ws.Range('E39').Value = '  -0.41%  '
$ws.Range('D40').Value = '7.259'
$ws.Range('E40').Value = '  +2.14%  '
$ws.Range('D41').Value = '2.871'
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('D42').Value = '0.5164'
$ws.Range('E42').Value = '  +2.17%  '
$ws.Range('D43').Value = '0.1697'
$ws.Range('E43').Value = '  +1.21%  '
$ws.Range('D44').Value = '8.535'
$ws.Range('E44').Value = '  -1.13%  '
$ws.Range('D45').Value = '10.64'
$ws.Range('E45').Value = '  +1.13%  '
$ws.Range('D46').Value = '0.4823'
$ws.Range('E46').Value = '  +2.35%  '
$ws.Range('D47').Value = '106.02'
$ws.Range('D48').Value = '1.932'
$ws.Range('E48').Value = '  +6.56%  '
$ws.Range('D49').Value = '0.9999'
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('D50').Value = '0.06336'
$ws.Range('E50').Value = '  -1.29%  '
$ws.Range('D51').Value = '1.659'
$ws.Range('E51').Value = '  -0.27%  '

# Remove the temporary text-format styling so cell styles match the original (no explicit style)
$ws.Range("D2:D51").ClearFormats()

